$d = $word.ActiveDocument

# Locate the paragraph right after "OUR PRODUCT: ..." - it is currently the
# existing empty paragraph that precedes "Stakeholders would include:".
# We insert 5 new paragraphs immediately before it:
#   1) an empty paragraph
#   2) "Why should I use this over the average tutoring site?"
#   3) "It is specifically designed for UNCC students and the courses you are taking."
#   4) "How will I be able to connect with a tutor that suits my needs?"
#   5) "You are able to filter through the available tutors by course, topic, and rating."

$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "OUR PRODUCT:*") {
        $anchor = $i
    }
}

$targetIndex = $anchor + 1

for ($i = 0; $i -lt 5; $i++) {
    $target = $d.Paragraphs.Item($targetIndex)
    $target.Range.InsertParagraphBefore()
}

$texts = @(
    "Why should I use this over the average tutoring site?",
    "It is specifically designed for UNCC students and the courses you are taking.",
    "How will I be able to connect with a tutor that suits my needs?",
    "You are able to filter through the available tutors by course, topic, and rating."
)

for ($i = 0; $i -lt 4; $i++) {
    $p = $d.Paragraphs.Item($targetIndex + 1 + $i)
    $r = $p.Range
    $r.End = $r.End - 1
    $r.Font.Size = 12
    $r.Font.SizeBi = 12
    $r.Font.NameAscii = "Times New Roman"
    $r.Font.NameFarEast = "Times New Roman"
    $r.Font.NameOther = "Times New Roman"
    $r.Font.NameBi = "Times New Roman"
    $r.Text = $texts[$i]
}
